# Fix bug: rename system codes and add new rows for "novo_sistema" system
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("systems")

# Update existing codes
$ws.Range("A3").Value = "sia-estacio"
$ws.Range("A4").Value = "financeiro"

# Add two new rows
$ws.Range("A5").Value = "novo_sistema"
$ws.Range("B5").Value = "Renomeado"

$ws.Range("A6").Value = "novo_sistem"
$ws.Range("B6").Value = "Novo Sistema"
